$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (row 2 through 23) from 45208 to 45212
$ws.Range("C2:C23").Value = 45212
